$d = $word.ActiveDocument
$para = $d.Paragraphs.Last
$r = $para.Range
$r.Collapse(0)
$r.InsertAfter("，使用给git创造分支")
$r.Font.NameAscii = "宋体"
$r.Font.NameFarEast = "宋体"
$r.Font.Name = "宋体"
$r.Font.NameBi = "宋体"
$r.Font.NameOther = "宋体"
$r.Font.Size = 14
